# Auto-generated Excel COM-interop edit script
# Updates H:N (price/profit) columns across all 8 item-category sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match the scheduled-runner refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2326.9524
$ws.Range("I51").Value = 2058.75
$ws.Range("J51").Value = 2492
$ws.Range("K51").Value = 2058.75
$ws.Range("L51").Value = 2492
$ws.Range("M51").Value = -1574.75
$ws.Range("N51").Value = -3460
$ws.Range("H127").Value = 2028.409
$ws.Range("I127").Value = 558.3125
$ws.Range("K127").Value = 1674.9375
$ws.Range("M127").Value = 3285.0625
$ws.Range("H129").Value = 904381.1
$ws.Range("I129").Value = 270.08334
$ws.Range("J129").Value = 1278496.1
$ws.Range("K129").Value = 810.2500200000001
$ws.Range("L129").Value = 3835488.3
$ws.Range("M129").Value = 4189.74998
$ws.Range("N129").Value = -3845488.3
$ws.Range("H131").Value = 1738
$ws.Range("J131").Value = 3604
$ws.Range("L131").Value = 10812
$ws.Range("N131").Value = -20892
$ws.Range("H132").Value = 1854.2646
$ws.Range("I132").Value = 1913.9062
$ws.Range("J132").Value = 900
$ws.Range("K132").Value = 5741.7186
$ws.Range("L132").Value = 2700
$ws.Range("M132").Value = -3211.7186
$ws.Range("N132").Value = -7760
$ws.Range("H137").Value = 861.93335
$ws.Range("I137").Value = 805.25
$ws.Range("K137").Value = 2415.75
$ws.Range("M137").Value = 134.25
$ws.Range("H138").Value = 1348.76
$ws.Range("I138").Value = 676.14514
$ws.Range("J138").Value = 2446.1843
$ws.Range("K138").Value = 2028.43542
$ws.Range("L138").Value = 7338.5529
$ws.Range("M138").Value = 3111.56458
$ws.Range("N138").Value = -17618.5529

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4825.49
$ws.Range("I32").Value = 3503.4285
$ws.Range("J32").Value = 9251.521000000001
$ws.Range("K32").Value = 3503.4285
$ws.Range("L32").Value = 9251.521000000001
$ws.Range("M32").Value = -3216.4285
$ws.Range("N32").Value = -9825.521000000001
$ws.Range("H45").Value = 1200.0667
$ws.Range("I45").Value = 966.8333
$ws.Range("J45").Value = 2133
$ws.Range("K45").Value = 966.8333
$ws.Range("L45").Value = 2133
$ws.Range("M45").Value = -589.8333
$ws.Range("N45").Value = -2887
$ws.Range("H110").Value = 945
$ws.Range("I110").Value = 795
$ws.Range("J110").Value = 1020
$ws.Range("K110").Value = 795
$ws.Range("L110").Value = 1020
$ws.Range("M110").Value = 1250
$ws.Range("N110").Value = -5110
$ws.Range("H132").Value = 1693.8
$ws.Range("I132").Value = 1256
$ws.Range("J132").Value = 2569.4
$ws.Range("K132").Value = 3768
$ws.Range("L132").Value = 7708.200000000001
$ws.Range("M132").Value = -1238
$ws.Range("N132").Value = -12768.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1196.7727
$ws.Range("I94").Value = 1238.6
$ws.Range("J94").Value = 1107.1428
$ws.Range("K94").Value = 1238.6
$ws.Range("L94").Value = 1107.1428
$ws.Range("M94").Value = -787.5999999999999
$ws.Range("N94").Value = -2009.1428
$ws.Range("H134").Value = 18427.25
$ws.Range("I134").Value = 1361.9592
$ws.Range("J134").Value = 94445.37
$ws.Range("K134").Value = 4085.8776
$ws.Range("L134").Value = 283336.11
$ws.Range("M134").Value = -1550.8776
$ws.Range("N134").Value = -288406.11

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2473.0784
$ws.Range("I31").Value = 2422.3784
$ws.Range("J31").Value = 2607.0715
$ws.Range("K31").Value = 2422.3784
$ws.Range("L31").Value = 2607.0715
$ws.Range("M31").Value = -2127.3784
$ws.Range("N31").Value = -3197.0715
$ws.Range("H34").Value = 2473.0784
$ws.Range("I34").Value = 2422.3784
$ws.Range("J34").Value = 2607.0715
$ws.Range("K34").Value = 2422.3784
$ws.Range("L34").Value = 2607.0715
$ws.Range("M34").Value = -2220.3784
$ws.Range("N34").Value = -3011.0715
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()
$ws.Range("H109").Value = 19145
$ws.Range("J109").Value = 19145
$ws.Range("L109").Value = 19145
$ws.Range("N109").Value = -21225
$ws.Range("H134").Value = 1104.6615
$ws.Range("I134").Value = 988.5738
$ws.Range("K134").Value = 2965.7214
$ws.Range("M134").Value = -430.7213999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 77742.8
$ws.Range("J37").Value = 77742.8
$ws.Range("L37").Value = 233228.4
$ws.Range("N37").Value = -233452.4
$ws.Range("H80").Value = 5250
$ws.Range("J80").Value = 7666.6665
$ws.Range("L80").Value = 22999.9995
$ws.Range("N80").Value = -24871.9995
$ws.Range("H83").Value = 5250
$ws.Range("J83").Value = 7666.6665
$ws.Range("L83").Value = 68999.9985
$ws.Range("N83").Value = -78359.9985
$ws.Range("H107").Value = 865017.7
$ws.Range("I107").Value = 1178.4
$ws.Range("J107").Value = 1944816.8
$ws.Range("K107").Value = 3535.2
$ws.Range("L107").Value = 5834450.4
$ws.Range("M107").Value = -1615.2
$ws.Range("N107").Value = -5838290.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H87").Value = 49900
$ws.Range("I87").Value = 50000
$ws.Range("J87").Value = 49800
$ws.Range("K87").Value = 50000
$ws.Range("L87").Value = 49800
$ws.Range("M87").Value = -48752
$ws.Range("N87").Value = -52296
$ws.Range("H90").Value = 49900
$ws.Range("I90").Value = 50000
$ws.Range("J90").Value = 49800
$ws.Range("K90").Value = 150000
$ws.Range("L90").Value = 149400
$ws.Range("M90").Value = -143760
$ws.Range("N90").Value = -161880
$ws.Range("H132").Value = 2881.44
$ws.Range("I132").Value = 3081.4666
$ws.Range("J132").Value = 2581.4
$ws.Range("K132").Value = 9244.399800000001
$ws.Range("L132").Value = 7744.200000000001
$ws.Range("M132").Value = -6714.399800000001
$ws.Range("N132").Value = -12804.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 222.03334
$ws.Range("I55").Value = 170.64706
$ws.Range("J55").Value = 289.23077
$ws.Range("K55").Value = 170.64706
$ws.Range("L55").Value = 289.23077
$ws.Range("M55").Value = 2.35293999999999
$ws.Range("N55").Value = -635.23077
$ws.Range("H58").Value = 990
$ws.Range("I58").Value = 990
$ws.Range("K58").Value = 990
$ws.Range("M58").Value = -730
$ws.Range("H100").Value = 2177.7778
$ws.Range("I100").Value = 2242.8572
$ws.Range("J100").Value = 1950
$ws.Range("K100").Value = 2242.8572
$ws.Range("L100").Value = 1950
$ws.Range("M100").Value = -1701.8572
$ws.Range("N100").Value = -3032

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 722.7
$ws.Range("I81").Value = 722.7
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1445.4
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -384.4000000000001
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 722.7
$ws.Range("I84").Value = 722.7
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 7227
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -1923
$ws.Range("N84").ClearContents()
$ws.Range("H93").Value = 26955.6
$ws.Range("J93").Value = 26955.6
$ws.Range("L93").Value = 26955.6
$ws.Range("N93").Value = -31947.6
$ws.Range("H96").Value = 4621
$ws.Range("I96").Value = 3933
$ws.Range("J96").Value = 5210.7144
$ws.Range("K96").Value = 3933
$ws.Range("L96").Value = 5210.7144
$ws.Range("M96").Value = -2560
$ws.Range("N96").Value = -7956.7144
$ws.Range("H123").Value = 45000
$ws.Range("J123").Value = 45000
$ws.Range("L123").Value = 45000
$ws.Range("N123").Value = -54800
